$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format before writing so that values
# like "1.001" / "236.23" / "9.010" are stored verbatim as strings instead
# of being auto-coerced into numbers (which would drop trailing zeros /
# change formatting). We restore the original (default) style afterwards
# so the cells keep the same "no explicit style" appearance as before.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.249.53"
$ws.Range("E2").Value = "  -1.32%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.859.62"
$ws.Range("E3").Value = "  -0.11%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").Value = "236.23"
$ws.Range("E5").Value = "  +0.55%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.13%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4655"
$ws.Range("E7").Value = "  -1.25%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2816"
$ws.Range("E8").Value = "  +2.31%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06393"
$ws.Range("E9").Value = "  +0.55%  "

# Row 10 - Solana
$ws.Range("D10").Value = "18.35"
$ws.Range("E10").Value = "  +4.29%  "

# Row 11 - becomes Litecoin (was WrappedEther)
$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").Value = "95.99"
$ws.Range("E11").Value = "  +12.85%  "

# Row 12 - becomes WrappedEther (was TRON)
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.852.86"
$ws.Range("E12").Value = "  -0.37%  "

# Row 13 - becomes TRON (was Litecoin)
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07520"
$ws.Range("E13").Value = "  +1.02%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "4.972"
$ws.Range("E14").Value = "  -4.81%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.6404"

# Row 16 - BitcoinCash
$ws.Range("D16").Value = "295.81"
$ws.Range("E16").Value = "  +21.44%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "30.246.76"
$ws.Range("E17").Value = "  -1.25%  "

# Row 18 - Dai
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.06%  "

# Row 19 - Avalanche
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.000007407"
$ws.Range("E20").Value = "  +0.52%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.108.90"
$ws.Range("E21").Value = "  +0.11%  "

# Row 22 - BinanceUSD
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.002"
$ws.Range("E23").Value = "  +0.05%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "6.043"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25 - Monero
$ws.Range("D25").Value = "164.85"
$ws.Range("E25").Value = "  -0.11%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "9.051"
$ws.Range("E26").Value = "  -3.16%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "19.35"
$ws.Range("E27").Value = "  +6.77%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "1.934"
$ws.Range("E28").Value = "  +2.31%  "

# Row 29 - Stellar
$ws.Range("D29").Value = "0.1083"
$ws.Range("E29").Value = "  +6.46%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "1.333"
$ws.Range("E30").Value = "  -3.53%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "4.035"
$ws.Range("E31").Value = "  -0.83%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.800"
$ws.Range("E32").Value = "  -1.78%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.04917"
$ws.Range("E33").Value = "  -0.18%  "

# Row 34 - ImmutableX
$ws.Range("D34").Value = "0.7254"
$ws.Range("E34").Value = "  +2.68%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "1.116"
$ws.Range("E35").Value = "  -3.07%  "

# Row 36 - HuobiToken
$ws.Range("D36").Value = "2.706"
$ws.Range("E36").Value = "  -0.19%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.01927"
$ws.Range("E37").Value = "  +0.61%  "

# Row 38 - MXToken
$ws.Range("D38").Value = "2.684"
$ws.Range("E38").Value = "  -0.06%  "

# Row 39 - becomes TrustWalletToken (was RenderToken)
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "0.8650"
$ws.Range("E39").Value = "  -1.74%  "

# Row 40 - becomes RenderToken (was TrustWalletToken)
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "1.969"
$ws.Range("E40").Value = "  -1.50%  "

# Row 41 - Quant
$ws.Range("D41").Value = "105.73"
$ws.Range("E41").Value = "  +0.36%  "

# Row 42 - PaxDollar
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.35%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "5.599"
$ws.Range("E43").Value = "  +0.95%  "

# Row 44 - TheSandbox
$ws.Range("D44").Value = "0.4063"
$ws.Range("E44").Value = "  -0.44%  "

# Row 45 - Aave
$ws.Range("D45").Value = "65.34"
$ws.Range("E45").Value = "  +3.11%  "

# Row 46 - Aptos
$ws.Range("D46").Value = "7.072"
$ws.Range("E46").Value = "  -2.84%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "9.010"
$ws.Range("E47").Value = "  +5.01%  "

# Row 48 - Algorand
$ws.Range("D48").Value = "0.1190"

# Row 49 - Elrond
$ws.Range("D49").Value = "33.85"
$ws.Range("E49").Value = "  +0.77%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.05565"
$ws.Range("E50").Value = "  +0.34%  "

# Row 51 - Decentraland
$ws.Range("D51").Value = "0.3732"
$ws.Range("E51").Value = "  +0.90%  "

# Restore the original (default/"Normal") style on column D so cells that
# had no explicit style keep that state in the saved file.
$ws.Range("D2:D51").Style = "Normal"
